# Updated cryptos list refresh (Price / Volume(1h) columns) on the "cryptos"
# worksheet. Only columns D (Price) and E (Volume(1h)) change; Coin / Link /
# index columns (A-C) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per data row (sheet row number matches the spreadsheet row).
# $null means that column is unchanged for that row in this refresh.
$updates = @(
    @{ Row = 2;  D = "26.885.88";    E = "  -0.59%  " }
    @{ Row = 3;  D = "1.859.88";     E = "  -0.11%  " }
    @{ Row = 4;  D = "0.9997";       E = "  -0.03%  " }
    @{ Row = 5;  D = "304.53";       E = "  -0.60%  " }
    @{ Row = 6;  D = "0.9997";       E = "  -0.05%  " }
    @{ Row = 7;  D = "0.5034";       E = "  -1.57%  " }
    @{ Row = 8;  D = $null;          E = "  -2.54%  " }
    @{ Row = 9;  D = "0.07167";      E = "  +0.71%  " }
    @{ Row = 10; D = "0.8931";       E = "  +0.16%  " }
    @{ Row = 11; D = "20.69";        E = "  +0.49%  " }
    @{ Row = 12; D = "1.875.84";     E = "  +0.75%  " }
    @{ Row = 13; D = "0.07483";      E = "  -0.79%  " }
    @{ Row = 14; D = "93.88";        E = "  +5.70%  " }
    @{ Row = 15; D = "5.227";        E = "  -1.35%  " }
    @{ Row = 16; D = "1.000";        E = "  +0.04%  " }
    @{ Row = 17; D = "0.000008491";  E = "  +1.35%  " }
    @{ Row = 18; D = "14.19";        E = "  +0.73%  " }
    @{ Row = 19; D = "0.9999";       E = "  -0.10%  " }
    @{ Row = 20; D = "26.930.38";    E = "  -0.56%  " }
    @{ Row = 21; D = "5.022";        E = "  -0.81%  " }
    @{ Row = 22; D = "2.111.06";     E = "  +0.80%  " }
    @{ Row = 23; D = $null;          E = "  -1.49%  " }
    @{ Row = 24; D = "6.412";        E = "  -0.93%  " }
    @{ Row = 25; D = "147.65";       E = "  -1.07%  " }
    @{ Row = 26; D = "1.775";        E = "  -3.61%  " }
    @{ Row = 27; D = "17.86";        E = "  -0.64%  " }
    @{ Row = 28; D = "2.084";        E = "  -0.28%  " }
    @{ Row = 29; D = "113.01";       E = "  +0.10%  " }
    @{ Row = 30; D = "4.688";        E = "  +0.03%  " }
    @{ Row = 31; D = "4.667";        E = "  +0.27%  " }
    @{ Row = 32; D = "0.09213";      E = "  +1.96%  " }
    @{ Row = 33; D = "0.05144";      E = "  +0.57%  " }
    @{ Row = 34; D = "0.7469";       E = "  +2.11%  " }
    @{ Row = 35; D = "2.967";        E = "  -2.93%  " }
    @{ Row = 36; D = $null;          E = "  -0.53%  " }
    @{ Row = 37; D = "3.250";        E = "  +6.46%  " }
    @{ Row = 38; D = "2.573";        E = "  +2.63%  " }
    @{ Row = 39; D = "0.02004";      E = "  -2.19%  " }
    @{ Row = 40; D = "0.5557";       E = "  +4.18%  " }
    @{ Row = 41; D = $null;          E = "  -0.20%  " }
    @{ Row = 42; D = "6.550";        E = "  -0.76%  " }
    @{ Row = 43; D = "117.37";       E = "  +1.77%  " }
    @{ Row = 44; D = "8.519";        E = "  +2.36%  " }
    @{ Row = 45; D = "0.1470";       E = $null }
    @{ Row = 46; D = "0.4676";       E = "  +1.14%  " }
    @{ Row = 47; D = "0.9993";       E = "  -0.07%  " }
    @{ Row = 48; D = "10.00";        E = "  +0.10%  " }
    @{ Row = 49; D = "1.561";        E = "  -0.25%  " }
    @{ Row = 50; D = $null;          E = "  -0.32%  " }
    @{ Row = 51; D = "62.93";        E = "  -1.89%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        # The "Price" column is stored as plain text in the sheet (it is
        # also used for thousands-dotted values like "26.885.88" that are
        # not valid numbers at all). Guard plain numeric-looking strings
        # with a leading apostrophe so Excel keeps storing them as text
        # instead of silently converting them to a Number cell, then drop
        # the resulting "quote prefix" formatting so the cell style is
        # left exactly as it was.
        if ($u.D -match '^[+-]?[0-9]*\.?[0-9]+$') {
            $cell.Value = "'" + $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
